$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so that values like
# "607.59" or trailing-zero values like "383.80" are preserved exactly,
# matching the inline-string cell type used in the source workbook,
# instead of being auto-converted to floating point numbers by Excel.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "72.929.57"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +5.85%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.655.77"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +6.34%  "
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "607.59"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.72%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "181.04"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.46%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +2.93%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.173"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +14.25%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "2.652.84"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +6.27%  "
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("E12").Value = "  +5.06%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.11"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("E14").Value = "  +10.49%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.134.95"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +6.03%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "27.02"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +5.66%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "72.845.31"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +5.89%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.654.81"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +5.86%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "383.80"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +6.81%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "11.60"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +7.07%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "7.93"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +5.68%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.24"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +5.61%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "2.03"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +23.00%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "73.52"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +5.38%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "4.46"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +7.23%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +12.23%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.792.94"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +6.31%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.43%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.0₃0976"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +11.36%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "542.27"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +7.18%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.11"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +5.49%  "
$ws.Range("E33").Value = "  +11.91%  "
$ws.Range("E34").Value = "  +4.82%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.16%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "162.47"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.17%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "19.40"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +4.45%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.42"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +9.65%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.113"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.08%  "
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("E41").Value = "  +10.36%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.67"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +17.20%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "5.14"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +8.55%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("E45").Value = "  +6.18%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "39.79"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.90%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "152.09"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.24%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "3.70"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +4.96%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.547"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +7.47%  "
$ws.Range("E50").Value = "  +10.98%  "
$ws.Range("E51").Value = "  +11.58%  "
